$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-1.3475791894612454, 1.4254374950085194, -0.268498227096124, 3.437979507283086)"
$ws.Range("C2").Value = "JSU(-1.2083159178476746, 0.975681428223216, 2.8081563777162395, 3.529756408571208)"
$ws.Range("D2").Value = "JSU(-0.8622460627788102, 1.135874730674988, 0.5215392680264201, 2.7304826363852626)"
$ws.Range("E2").Value = "NIG(0.7491766071120747, 0.5089096743473225, 5.705237130142033, 5.129588123281961)"
